# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (column G) holds re-simulated values for each saved game
# row; this writes the freshly calculated s_vals into column G for every
# row on the sheet (rows 2-60), leaving every other column untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 0
    4  = 0
    5  = 0
    6  = 0
    7  = 0
    8  = 1
    9  = 2
    10 = 3
    11 = 1
    12 = 1
    13 = 3
    14 = 0
    16 = 1
    17 = 0
    18 = 2
    19 = 0
    20 = 4
    22 = 0
    23 = 2
    24 = 1
    25 = 0
    26 = 0
    27 = 1
    28 = 1
    29 = 3
    30 = 0
    31 = 2
    32 = 0
    33 = 1
    34 = 0
    35 = 1
    36 = 1
    37 = 2
    38 = 2
    39 = 0
    40 = 0
    41 = 1
    42 = 0
    43 = 1
    44 = 2
    45 = 0
    46 = 0
    47 = 2
    48 = 0
    49 = 2
    50 = 0
    51 = 0
    52 = 0
    53 = 1
    55 = 2
    56 = 1
    57 = 2
    59 = 1
    60 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
